$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 (I0) and J1 (IF), matching the style of existing headers ---
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-70 for columns I (I0) and J (IF) ---
$rowData = @(
    @{Row=2; I=1; J=1},
    @{Row=3; I=8; J=8},
    @{Row=4; I=7; J=8},
    @{Row=5; I=6; J=6},
    @{Row=6; I=8; J=9},
    @{Row=7; I=9; J=9},
    @{Row=8; I=8; J=9},
    @{Row=9; I=1; J=2},
    @{Row=10; I=8; J=8},
    @{Row=11; I=7; J=7},
    @{Row=12; I=7; J=7},
    @{Row=13; I=7; J=7},
    @{Row=14; I=9; J=9},
    @{Row=15; I=7; J=9},
    @{Row=16; I=6; J=6},
    @{Row=17; I=9; J=9},
    @{Row=18; I=9; J=9},
    @{Row=19; I=5; J=6},
    @{Row=20; I=9; J=9},
    @{Row=21; I=8; J=8},
    @{Row=22; I=8; J=8},
    @{Row=23; I=7; J=7},
    @{Row=24; I=7; J=7},
    @{Row=25; I=9; J=9},
    @{Row=26; I=7; J=7},
    @{Row=27; I=6; J=7},
    @{Row=28; I=6; J=6},
    @{Row=29; I=7; J=7},
    @{Row=30; I=5; J=5},
    @{Row=31; I=9; J=9},
    @{Row=32; I=5; J=6},
    @{Row=33; I=9; J=9},
    @{Row=34; I=8; J=8},
    @{Row=35; I=10; J=10},
    @{Row=36; I=8; J=8},
    @{Row=37; I=5; J=5},
    @{Row=38; I=7; J=7},
    @{Row=39; I=5; J=5},
    @{Row=40; I=6; J=7},
    @{Row=41; I=3; J=4},
    @{Row=42; I=7; J=7},
    @{Row=43; I=8; J=8},
    @{Row=44; I=4; J=5},
    @{Row=45; I=7; J=8},
    @{Row=46; I=9; J=9},
    @{Row=47; I=6; J=6},
    @{Row=48; I=6; J=6},
    @{Row=49; I=10; J=10},
    @{Row=50; I=5; J=6},
    @{Row=51; I=8; J=8},
    @{Row=52; I=8; J=8},
    @{Row=53; I=8; J=8},
    @{Row=54; I=6; J=6},
    @{Row=55; I=5; J=5},
    @{Row=56; I=7; J=7},
    @{Row=57; I=5; J=6},
    @{Row=58; I=8; J=8},
    @{Row=59; I=8; J=8},
    @{Row=60; I=7; J=7},
    @{Row=61; I=7; J=7},
    @{Row=62; I=6; J=6},
    @{Row=63; I=7; J=7},
    @{Row=64; I=7; J=7},
    @{Row=65; I=4; J=4},
    @{Row=66; I=6; J=7},
    @{Row=67; I=5; J=5},
    @{Row=68; I=7; J=7},
    @{Row=69; I=6; J=6},
    @{Row=70; I=5; J=5}
)

foreach ($item in $rowData) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
}

Write-Output "Updated $($rowData.Count) rows with I0/IF columns"
